$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1633.3334
$ws.Range("J121").Value = 1450
$ws.Range("L121").Value = 4350
$ws.Range("N121").Value = -7844

$ws.Range("H124").Value = 22763.25
$ws.Range("J124").Value = 22763.25
$ws.Range("L124").Value = 22763.25
$ws.Range("N124").Value = -32583.25

$ws.Range("H132").Value = 5440067
$ws.Range("I132").Value = 6762670.5
$ws.Range("J132").Value = 2698.4443
$ws.Range("K132").Value = 20288011.5
$ws.Range("L132").Value = 8095.3329
$ws.Range("M132").Value = -20285481.5
$ws.Range("N132").Value = -13155.3329

$ws.Range("H141").Value = 2810.5625
$ws.Range("I141").Value = 2289.9167
$ws.Range("K141").Value = 6869.750100000001
$ws.Range("M141").Value = -1689.750100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 10455.25
$ws.Range("J98").Value = 10455.25
$ws.Range("L98").Value = 10455.25
$ws.Range("N98").Value = -16445.25

$ws.Range("H122").Value = 2481.524
$ws.Range("I122").Value = 2591.111
$ws.Range("J122").Value = 2399.3333
$ws.Range("K122").Value = 7773.333
$ws.Range("L122").Value = 7197.999899999999
$ws.Range("M122").Value = -5323.333
$ws.Range("N122").Value = -12097.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 150
$ws.Range("K2").Value = 150
$ws.Range("M2").Value = -37

$ws.Range("H16").Value = 750.6667
$ws.Range("I16").Value = 570.3333
$ws.Range("J16").Value = 931
$ws.Range("K16").Value = 570.3333
$ws.Range("L16").Value = 931
$ws.Range("M16").Value = -283.3333
$ws.Range("N16").Value = -1505

$ws.Range("H31").Value = 26041.477
$ws.Range("I31").Value = 39931.883
$ws.Range("J31").Value = 3469.5625
$ws.Range("K31").Value = 39931.883
$ws.Range("L31").Value = 3469.5625
$ws.Range("M31").Value = -39636.883
$ws.Range("N31").Value = -4059.5625

$ws.Range("H34").Value = 26041.477
$ws.Range("I34").Value = 39931.883
$ws.Range("J34").Value = 3469.5625
$ws.Range("K34").Value = 39931.883
$ws.Range("L34").Value = 3469.5625
$ws.Range("M34").Value = -39729.883
$ws.Range("N34").Value = -3873.5625

$ws.Range("H68").Value = 14449.583
$ws.Range("J68").Value = 14449.583
$ws.Range("L68").Value = 14449.583
$ws.Range("N68").Value = -15947.583

$ws.Range("H71").Value = 14449.583
$ws.Range("J71").Value = 14449.583
$ws.Range("L71").Value = 43348.749
$ws.Range("N71").Value = -50836.749

$ws.Range("H74").Value = 23540.428
$ws.Range("J74").Value = 23540.428
$ws.Range("L74").Value = 23540.428
$ws.Range("N74").Value = -25288.428

$ws.Range("H77").Value = 23540.428
$ws.Range("J77").Value = 23540.428
$ws.Range("L77").Value = 70621.284
$ws.Range("N77").Value = -79357.284

$ws.Range("H99").Value = 19665.666
$ws.Range("I99").Value = 3993.3333
$ws.Range("J99").Value = 35338
$ws.Range("K99").Value = 3993.3333
$ws.Range("L99").Value = 35338
$ws.Range("M99").Value = -2495.3333
$ws.Range("N99").Value = -38334

$ws.Range("H107").Value = 768.0625
$ws.Range("I107").Value = 967.8889
$ws.Range("J107").Value = 511.14285
$ws.Range("K107").Value = 967.8889
$ws.Range("L107").Value = 511.14285
$ws.Range("M107").Value = 952.1111
$ws.Range("N107").Value = -4351.14285

$ws.Range("H113").Value = 750.6667
$ws.Range("I113").Value = 570.3333
$ws.Range("J113").Value = 931
$ws.Range("K113").Value = 570.3333
$ws.Range("L113").Value = 931
$ws.Range("M113").Value = 1599.6667
$ws.Range("N113").Value = -5271

$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 675
$ws.Range("K122").Value = 2025
$ws.Range("M122").Value = 425

$ws.Range("H126").Value = 19665.666
$ws.Range("I126").Value = 3993.3333
$ws.Range("J126").Value = 35338
$ws.Range("K126").Value = 11979.9999
$ws.Range("L126").Value = 106014
$ws.Range("M126").Value = -9509.999899999999
$ws.Range("N126").Value = -110954

$ws.Range("H141").Value = 65981.82
$ws.Range("J141").Value = 34475
$ws.Range("L141").Value = 34475
$ws.Range("N141").Value = -44835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 13948.5625
$ws.Range("I5").Value = 1755.6
$ws.Range("J5").Value = 19490.818
$ws.Range("K5").Value = 5266.799999999999
$ws.Range("L5").Value = 58472.454
$ws.Range("M5").Value = -5154.799999999999
$ws.Range("N5").Value = -58696.454

$ws.Range("H12").Value = 33
$ws.Range("J12").Value = 36.235294
$ws.Range("L12").Value = 108.705882
$ws.Range("N12").Value = -454.705882

$ws.Range("H34").Value = 1360

$ws.Range("H38").Value = 109.8
$ws.Range("I38").Value = 11.25
$ws.Range("J38").Value = 222.42857
$ws.Range("K38").Value = 33.75
$ws.Range("L38").Value = 667.28571
$ws.Range("M38").Value = 313.25
$ws.Range("N38").Value = -1361.28571

$ws.Range("H39").Value = 18750
$ws.Range("J39").Value = 18750
$ws.Range("L39").Value = 56250
$ws.Range("N39").Value = -56838

$ws.Range("H55").Value = 11971.526
$ws.Range("J55").Value = 7462.294
$ws.Range("L55").Value = 22386.882
$ws.Range("N55").Value = -22740.882

$ws.Range("H131").Value = 825.47
$ws.Range("J131").Value = 867.0227
$ws.Range("L131").Value = 2601.0681
$ws.Range("N131").Value = -12681.0681

$ws.Range("H135").Value = 13948.5625
$ws.Range("I135").Value = 1755.6
$ws.Range("J135").Value = 19490.818
$ws.Range("K135").Value = 15800.4
$ws.Range("L135").Value = 175417.362
$ws.Range("M135").Value = -13265.4
$ws.Range("N135").Value = -180487.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H64").Value = 48125.5
$ws.Range("J64").Value = 48125.5
$ws.Range("L64").Value = 48125.5
$ws.Range("N64").Value = -48621.5

$ws.Range("H67").Value = 48125.5
$ws.Range("J67").Value = 48125.5
$ws.Range("L67").Value = 48125.5
$ws.Range("N67").Value = -49841.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3498.75
$ws.Range("I46").Value = 501
$ws.Range("J46").Value = 3927
$ws.Range("K46").Value = 501
$ws.Range("L46").Value = 3927
$ws.Range("M46").Value = -313
$ws.Range("N46").Value = -4303

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 30445
$ws.Range("J124").Value = 30445
$ws.Range("L124").Value = 30445
$ws.Range("N124").Value = -40265

$ws.Range("H140").Value = 57484.785
$ws.Range("J140").Value = 57484.785
$ws.Range("L140").Value = 57484.785
$ws.Range("N140").Value = -67844.785

$ws.Range("H141").Value = 49310.875
$ws.Range("J141").Value = 49310.875
$ws.Range("L141").Value = 49310.875
$ws.Range("N141").Value = -59670.875
